$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: move the "tenure track" marker from column G to column J ---
$ws.Range("G14").ClearContents()
$ws.Range("J14").Value = "tenure track"

# --- Row 15: move the "tenure track" marker from column G to column J ---
$ws.Range("G15").ClearContents()
$ws.Range("J15").Value = "tenure track"

# --- Row 20: new posting for East Carolina University ---
$ws.Range("B20").Value = 17
$ws.Range("C20").Value = "East Carolina University"

$dueDate = Get-Date -Year 2012 -Month 11 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Range("D20").Value = $dueDate
$ws.Range("D20").NumberFormat = $ws.Range("E14").NumberFormat

$ws.Range("E20").Value = "https://ecu.peopleadmin.com/applicants/jsp/shared/position/JobDetails_css.jsp"
$ws.Range("F20").Value = "ASSISTANT/ASSOCIATE PROFESSOR  "
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = "A minimum of three years of managerial experience in the construction industry; evidence of successful teaching in construction management at the college/university level; demonstrated ability to secure external funding, grants, and/or industry support and professional certification. Preference will be given to candidates able to teach courses in the following areas: Building Information Modeling (BIM), Mechanical, Electrical, and Plumbing (MEP) Systems, Estimating, Scheduling, and Cost Control, in face-to-face and distance learning environments. "
$ws.Range("I20").Value = "An earned PhD degree in Construction Management, Civil Engineering, or closely related field and a demonstrated potential to develop and sustain an active research agenda, effective communication and interpersonal skills; ability and desire to work in a team setting and make positive contributions to the department. All degrees must be from a regionally accredited institution. ABD may be considered for appointment at the rank of instructor and must provide documented evidence that their dissertation defense and completion of study is completed prior to the fall 2017 semester for fall 2017 hire. "
$ws.Range("K20").Value = "Candidates must submit a cover letter, a curriculum vitae/resume, teaching statement, research statement and a list of three references, including contact information, online.  "

# Row 20 grew to accommodate the new text (matches the ~160pt height used by Excel)
$ws.Rows.Item(20).RowHeight = 160

# Reflect the new active selection after entering the data
$ws.Range("B20").Select()
